# Update "想去人数" (column F) values across the four sheets to match
# the regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 1137
$ws.Cells.Item(5, 6).Value = 54
$ws.Cells.Item(7, 6).Value = 785
$ws.Cells.Item(14, 6).Value = 938
$ws.Cells.Item(15, 6).Value = 126
$ws.Cells.Item(16, 6).Value = 2050
$ws.Cells.Item(17, 6).Value = 527
$ws.Cells.Item(18, 6).Value = 8709
$ws.Cells.Item(19, 6).Value = 804
$ws.Cells.Item(21, 6).Value = 78
$ws.Cells.Item(22, 6).Value = 104
$ws.Cells.Item(23, 6).Value = 26
$ws.Cells.Item(24, 6).Value = 237

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 526
$ws.Cells.Item(9, 6).Value = 129
$ws.Cells.Item(10, 6).Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 5595
$ws.Cells.Item(3, 6).Value = 429

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 5595
$ws.Cells.Item(4, 6).Value = 429
$ws.Cells.Item(6, 6).Value = 526
$ws.Cells.Item(7, 6).Value = 1137
$ws.Cells.Item(10, 6).Value = 54
$ws.Cells.Item(12, 6).Value = 785
$ws.Cells.Item(23, 6).Value = 938
$ws.Cells.Item(24, 6).Value = 126
$ws.Cells.Item(25, 6).Value = 129
$ws.Cells.Item(26, 6).Value = 3
$ws.Cells.Item(27, 6).Value = 2050
$ws.Cells.Item(28, 6).Value = 527
$ws.Cells.Item(29, 6).Value = 8709
$ws.Cells.Item(32, 6).Value = 804
$ws.Cells.Item(34, 6).Value = 78
$ws.Cells.Item(35, 6).Value = 104
$ws.Cells.Item(37, 6).Value = 26
$ws.Cells.Item(39, 6).Value = 237
